$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "No." column values (12 -> 1, 14 -> 2, 16 -> 3) ---
# Force the cells to stay as text (they already are shared strings), then
# restore the default "Normal" style so no extra number-format style sticks.
$noRange = $ws.Range("A2:A4")
$noRange.NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("A3").Value = "2"
$ws.Range("A4").Value = "3"
$noRange.Style = "Normal"

# --- Resize the three profile-image pictures (381000x190500 -> 476250x285750 EMU) ---
# 476250 EMU = 37.5 pt, 285750 EMU = 22.5 pt (1 pt = 12700 EMU)
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Width = 37.5
    $shp.Height = 22.5
}

# --- Narrow column B from 30 to 25 characters wide ---
$ws.Columns("B").ColumnWidth = 25 - (5/6)

# --- Give rows 2-4 a fixed custom height of 40 ---
$ws.Rows(2).RowHeight = 40
$ws.Rows(3).RowHeight = 40
$ws.Rows(4).RowHeight = 40
